# ---------------------------------------------------------------------------
# Applies three edits described by the commit's diff:
#   1. Insert a new paragraph ("Back to the main.SCSS file") right after the
#      "This is how you Nest" heading paragraph.
#   2. Append a new run (" Add the class of main__p to your paragraph tag.")
#      to the paragraph that ends with
#      "...have a class attached to the paragraph."
#   3. Mark the run holding the "Our Regular CSS code will look like this."
#      picture (anchorId 29695F51) as NoProof (adds <w:rPr><w:noProof/></w:rPr>).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Namespace + pkg wrapper helpers used for precise OOXML insertion through
# Range.InsertXML (which REPLACES the target range's contents).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $bodyXml + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# Edit 1: new paragraph after "This is how you Nest"
# ---------------------------------------------------------------------------

$headingRange = $d.Content
$found = $headingRange.Find.Execute("This is how you Nest", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find heading 'This is how you Nest'"
}
$headingRange.Collapse(0)
$null = $headingRange.InsertParagraphAfter()

# Locate the freshly inserted (still empty) paragraph - it immediately
# follows the heading paragraph.
$newParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "This is how you Nest`r") {
        $newParaIdx = $i + 1
        break
    }
}
if ($newParaIdx -eq -1) {
    throw "Could not locate newly inserted paragraph after heading"
}
$newPara = $d.Paragraphs($newParaIdx)
$newParaBody = '<w:p ' + $wNs + '>' + `
    '<w:r><w:t xml:space="preserve">Back to the </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>main.SCSS</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>file</w:t></w:r>' + `
    '</w:p>'
$null = $newPara.Range.InsertXML((New-PkgXml $newParaBody))

# ---------------------------------------------------------------------------
# Edit 2: append a run to the "...attached to the paragraph." paragraph
# ---------------------------------------------------------------------------

$classParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Remember that with a class we also have to change it in the html to have a class attached to the paragraph.*") {
        $classParaIdx = $i
        break
    }
}
if ($classParaIdx -eq -1) {
    throw "Could not find the 'Remember that with a class...' paragraph"
}
$classPara = $d.Paragraphs($classParaIdx)
$classRange = $classPara.Range

# Recover this paragraph's own attributes (w14:paraId, rsids, ...) so they
# survive the InsertXML-based rewrite.
$openTag = '<w:p>'
$openXml = $classRange.WordOpenXML
if ($openXml -match '(<w:p\b[^>]*>)') {
    $openTag = $matches[1]
}

$classBody = $openTag + `
    '<w:r><w:lastRenderedPageBreak/><w:t>Remember that with a class we also have to change it in the html to have a class attached to the paragraph.</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> Add the class of main__p to your paragraph tag.</w:t></w:r>' + `
    '</w:p>'
$null = $classRange.InsertXML((New-PkgXml $classBody))

# ---------------------------------------------------------------------------
# Edit 3: mark the "Our Regular CSS code will look like this." picture run
#         as NoProof
# ---------------------------------------------------------------------------

$cssParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Our Regular CSS code will look like this.`r") {
        $cssParaIdx = $i
        break
    }
}
if ($cssParaIdx -eq -1) {
    throw "Could not find 'Our Regular CSS code will look like this.' paragraph"
}
$picturePara = $d.Paragraphs($cssParaIdx + 1)
$picturePara.Range.NoProofing = $true

Write-Output "Done."
